$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: 01/01/2018 -> 01/01/2021 (row 8)
# Pre-format as text so Excel keeps the date-looking string as a literal
# value instead of silently converting it to a date serial number.
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# Docentes responsáveis: (row 13)
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# Método: (row 19)
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# Critério: (row 20)
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas."

# Norma de recuperação: (row 21)
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# Bibliografia: (row 22)
$bibliografia = "KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.KOTLER, P.; KARTAJAYA, H.; SETIAWAN, I. Marketing 4.0: do Tradicional ao Digital. São Paulo: Sextante, 2017.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L.  Marketing Essencial. 5 ed. São Paulo: Pearson, 2013.SANDHUSEN, R. L. Marketing Básico - Série Essencial. 3 ed. São Paulo: Saraiva, 2010.SAPIRO, Arão., CHIAVENATO, I. Planejamento Estratégico. Campus, 2ª. edição, 2010 KOTLER, P. Administração de Marketing, edição do milênio, revisão técnica de Prof. Arão Sapiro. Prentice-Hall, 2000. HOOLEY, Graham J.; PIERCY, Nigel F.; SAUNDERS, John A. Estratégia de Marketing e Posicionamento Competitivo tradução e revisão técnica: Prof. Arão Sapiro. Pearson Education do Brasil, 2001. SAPIRO, ARAO; GANGANA, MAURÍCIO; LIMA, MIGUEL; VILHENA, JOÃO BAPTISTA. Gestão de Marketing . FGV Editora, 2004. BOONE, L. e KURTZ, D.L. Marketing contemporâneo. 8ª ed. São Paulo, Livros Técnicos e Científicos, 1998. KOTLER, P; JATURISPITAK, S. e MAESINCIE, S. O marketing das nações. São Paulo, Futura, 1997. MARTINS, J.R. e BLECHER, N. O império das marcas. 2ª ed. São Paulo, Negócio Editora, 1997 THUROW, L.C. O futuro do capitalismo. 2ª ed. São Paulo, Rocco, 1997. VAZ, G. N. Marketing institucional. São Paulo, Pioneira, 1995. Bibliografia Complementar Artigos das Revistas: Marketing, Meio e Mensagem, Exame, Dinheiro, Revista da Escola de Administração da FEA-USP, Revista ESPM."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
